$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All touched cells in this sheet hold text values (coin names, URLs,
# price strings using "." as thousands separators, and padded percent
# strings). Force text number-format first so Excel does not silently
# reinterpret numeric-looking strings (e.g. "73.00", "0.0751") as numbers
# and strip formatting/precision.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.828.98'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +4.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.276.00'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.28%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.13'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.67'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +5.21%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.75%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.488'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.73'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.02'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +6.29%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.33%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.73%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.626.07'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.32'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.251.16'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.765'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.770.85'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.30'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +9.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0910'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.18%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.25'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '242.95'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +5.42%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.93'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.34'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.82%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.62'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.04%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.64'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '34.05'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +6.60%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0751'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.75%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Celestia'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.90'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +8.26%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.38'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.45%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +5.40%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.93'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.074.78'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.71'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.66%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.27%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.93'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +7.36%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +8.55%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.16'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.07%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.53'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.28%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.00'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +7.53%  '
